$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "2035_TM152_EIR_Alt1_06" row right after the existing
# "2035_TM152_EIR_Alt1_05" row (row 161). Copying the row and doing an
# "insert copied cells" at that same row (shift down) duplicates its
# formatting exactly the way Excel does: the original row keeps its place
# and becomes a plain interior row, while the pasted copy captures the
# bottom-of-group border and moves the surviving content down to row 162.
$ws.Rows.Item(161).Copy()
$ws.Rows.Item(161).Insert(-4121)
$ws.Range("C162").Value = "2035_TM152_EIR_Alt1_06"

# --- Insert the new "2050_TM152_EIR_Alt1_06" row right after the existing
# "2050_TM152_EIR_Alt1_05" row. After the insert above, that row is now 167.
$ws.Rows.Item(167).Copy()
$ws.Rows.Item(167).Insert(-4121)
$ws.Range("C168").Value = "2050_TM152_EIR_Alt1_06"

# --- Update the AutoFilter defined name range to reflect the two new rows.
$wb.Names("_xlnm._FilterDatabase").RefersToR1C1 = "=all_runs!R1C1:R175C8"

# --- Match the saved view state (frozen pane / selection / window size).
$ws.Application.ActiveWindow.ScrollRow = 134
$ws.Range("A168").Select()

Write-Host "done"
